# Rename sheet "nad" to "nadp"
$wb = $excel.ActiveWorkbook
$wsNad = $wb.Worksheets.Item("nad")
$wsNad.Name = "nadp"

# Select the "dna" worksheet and populate new columns F:H with data rows 2-25
$wsDna = $wb.Worksheets.Item("dna")

$data = @(
    @(45479, 47505, 47532),
    @(594246, 592748, 607592),
    @(1094433, 1084946, 1085169),
    @(2116238, 2083453, 2167345),
    @(3969176, 4062854, 4039748),
    @(7683769, 7525415, 7405485),
    @(14950687, 15095009, 15088696),
    @(32258614, 31554882, 32605324),
    @(5291536, 5901404, 4604966),
    @(5528206, 5739476, 5846970),
    @(5814092, 4987666, 5513712),
    @(4128609, 4672424, 6330959),
    @(5736942, 6176194, 4683002),
    @(4717836, 5089783, 5281606),
    @(1959179, 1769960, 2367081),
    @(2694377, 2093140, 2698278),
    @(4573632, 3935763, 4444847),
    @(4350526, 3753182, 5738998),
    @(5279350, 5743242, 8235940),
    @(3540809, 4081133, 4407561),
    @(2773509, 3468420, 4423282),
    @(1940053, 1828854, 2571315),
    @(1988144, 2075039, 2073654),
    @(1427410, 1669692, 1530632)
)

$row = 2
foreach ($rowData in $data) {
    $wsDna.Cells.Item($row, 6).Value = $rowData[0]
    $wsDna.Cells.Item($row, 7).Value = $rowData[1]
    $wsDna.Cells.Item($row, 8).Value = $rowData[2]
    $row++
}

# Update selection on dna sheet to F18:H25
$wsDna.Range("F18:H25").Select()

# Restore the "nadp" sheet as the active tab (matches original workbook state)
$wsNad.Activate()
